## #1509 fix doc, core jar delivery now to different directory.
## Update the "Delivered files and required copy actions" table:
##  - insert a new "delivered jars"/path column between "component" and
##    the existing jar-list column
##  - rename the old "build delivers in component/target/" header to
##    "build delivery directory"
##  - resize the table's columns
##  - fix "core" row to say the jars are delivered to component/build
##    (instead of component/target)

function Insert-CellBodyXml($cell, $bodyXml) {
    $rng = $cell.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      $bodyXml +
      '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

function Set-NewColumnCell($table, $row, $plainText) {
    $cell = $table.Cell($row, 2)
    $body = '<w:body><w:p><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:r><w:t>' + $plainText + '</w:t></w:r></w:p></w:body>'
    Insert-CellBodyXml $cell $body
    $t2 = $word.ActiveDocument.Tables.Item(2)
    $cell2 = $t2.Cell($row, 2)
    $cell2.Range.Paragraphs.Item(1).Range.Delete()
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# --- Insert the new 2nd column (before the current jars/description column) ---
$existingCol2 = $t.Columns.Item(2)
$t.Columns.Add($existingCol2) | Out-Null

# --- Re-fetch the table and set the final column widths (in points; values
#     below are the twentieths-of-a-point (dxa) target widths / 20) ---
$t = $d.Tables.Item(2)
$t.Columns.Item(1).Width = 89.6
$t.Columns.Item(2).Width = 104.05
$t.Columns.Item(3).Width = 123.3
$t.Columns.Item(4).Width = 144.55

# --- Header row ---
$t = $d.Tables.Item(2)
$headerDirCell = $t.Cell(1, 2)
$headerDirBody = '<w:body><w:p><w:pPr><w:keepNext/><w:keepLines/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>build delivery directory</w:t></w:r></w:p></w:body>'
Insert-CellBodyXml $headerDirCell $headerDirBody
$t2 = $d.Tables.Item(2)
$t2.Cell(1, 2).Range.Paragraphs.Item(1).Range.Delete()

$t = $d.Tables.Item(2)
$headerJarsCell = $t.Cell(1, 3)
$headerBody = '<w:body><w:p><w:pPr><w:keepNext/><w:keepLines/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>delivered jars</w:t></w:r></w:p></w:body>'
Insert-CellBodyXml $headerJarsCell $headerBody
$t2 = $d.Tables.Item(2)
$t2.Cell(1, 3).Range.Paragraphs.Item(1).Range.Delete()

# --- Data rows: new column 2 contents ---
$t = $d.Tables.Item(2)
Set-NewColumnCell $t 2 "/build"

$t = $d.Tables.Item(2)
Set-NewColumnCell $t 3 "/target"

$t = $d.Tables.Item(2)
Set-NewColumnCell $t 4 "/target"

$t = $d.Tables.Item(2)
Set-NewColumnCell $t 5 "/garget"
